$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Fn1"
$ws.Cells.Item(2, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 24.721787
$ws.Cells.Item(2, 8).Value = 74.16536099999999
$ws.Cells.Item(2, 9).Value = 0.01340847400407589
$ws.Cells.Item(2, 10).Value = 0.01345828544885583
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 2.311298
$ws.Cells.Item(2, 14).Value = 6.933894
$ws.Cells.Item(2, 15).Value = 0.8122763614007964
$ws.Cells.Item(2, 16).Value = 0.863868705896424
$ws.Cells.Item(2, 17).Value = 57.139416849526
$ws.Cells.Item(2, 18).Value = 514.254751645734
$ws.Cells.Item(2, 19).Value = 0.01089138647596793
$ws.Cells.Item(2, 20).Value = 0.01162619163428776

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Fn1"
$ws.Cells.Item(3, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(3, 4).Value = "M1"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 24.721787
$ws.Cells.Item(3, 8).Value = 74.16536099999999
$ws.Cells.Item(3, 9).Value = 0.01340847400407589
$ws.Cells.Item(3, 10).Value = 0.01345828544885583
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.02434666666666667
$ws.Cells.Item(3, 14).Value = 0.07304
$ws.Cells.Item(3, 15).Value = 0.008556327142686946
$ws.Cells.Item(3, 16).Value = 0.009099788701510982
$ws.Cells.Item(3, 17).Value = 0.6018931074933332
$ws.Cells.Item(3, 18).Value = 5.417037967439999
$ws.Cells.Item(3, 19).Value = 0.0001147272900630869
$ws.Cells.Item(3, 20).Value = 0.0001224675538692079

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Fn1"
$ws.Cells.Item(4, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 24.721787
$ws.Cells.Item(4, 8).Value = 74.16536099999999
$ws.Cells.Item(4, 9).Value = 0.01340847400407589
$ws.Cells.Item(4, 10).Value = 0.01345828544885583
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.509813
$ws.Cells.Item(4, 14).Value = 1.019626
$ws.Cells.Item(4, 15).Value = 0.1791673114565167
$ws.Cells.Item(4, 16).Value = 0.1270315054020651
$ws.Cells.Item(4, 17).Value = 12.603488395831
$ws.Cells.Item(4, 18).Value = 75.62093037498599
$ws.Cells.Item(4, 19).Value = 0.002402360238044873
$ws.Cells.Item(4, 20).Value = 0.001709626260698864

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Fn1"
$ws.Cells.Item(5, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(5, 4).Value = "FAPs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1458.280985666667
$ws.Cells.Item(5, 8).Value = 4374.842957
$ws.Cells.Item(5, 9).Value = 0.7909348416823457
$ws.Cells.Item(5, 10).Value = 0.7938731034993884
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 2.311298
$ws.Cells.Item(5, 14).Value = 6.933894
$ws.Cells.Item(5, 15).Value = 0.8122763614007964
$ws.Cells.Item(5, 16).Value = 0.863868705896424
$ws.Cells.Item(5, 17).Value = 3370.521925609396
$ws.Cells.Item(5, 18).Value = 30334.69733048456
$ws.Cells.Item(5, 19).Value = 0.6424576753068507
$ws.Cells.Item(5, 20).Value = 0.6858021305659945

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Fn1"
$ws.Cells.Item(6, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(6, 4).Value = "M1"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1458.280985666667
$ws.Cells.Item(6, 8).Value = 4374.842957
$ws.Cells.Item(6, 9).Value = 0.7909348416823457
$ws.Cells.Item(6, 10).Value = 0.7938731034993884
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.02434666666666667
$ws.Cells.Item(6, 14).Value = 0.07304
$ws.Cells.Item(6, 15).Value = 0.008556327142686946
$ws.Cells.Item(6, 16).Value = 0.009099788701510982
$ws.Cells.Item(6, 17).Value = 35.50428106436444
$ws.Cells.Item(6, 18).Value = 319.53852957928
$ws.Cells.Item(6, 19).Value = 0.006767497253983457
$ws.Cells.Item(6, 20).Value = 0.007224077497657194

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Fn1"
$ws.Cells.Item(7, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1458.280985666667
$ws.Cells.Item(7, 8).Value = 4374.842957
$ws.Cells.Item(7, 9).Value = 0.7909348416823457
$ws.Cells.Item(7, 10).Value = 0.7938731034993884
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.509813
$ws.Cells.Item(7, 14).Value = 1.019626
$ws.Cells.Item(7, 15).Value = 0.1791673114565167
$ws.Cells.Item(7, 16).Value = 0.1270315054020651
$ws.Cells.Item(7, 17).Value = 743.4506041456802
$ws.Cells.Item(7, 18).Value = 4460.703624874081
$ws.Cells.Item(7, 19).Value = 0.1417096691215116
$ws.Cells.Item(7, 20).Value = 0.1008468954357368

# Row 8
$ws.Cells.Item(8, 1).Value = "M1"
$ws.Cells.Item(8, 2).Value = "Fn1"
$ws.Cells.Item(8, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 60.09979000000001
$ws.Cells.Item(8, 8).Value = 180.29937
$ws.Cells.Item(8, 9).Value = 0.03259661091107292
$ws.Cells.Item(8, 10).Value = 0.03271770480169137
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 2.311298
$ws.Cells.Item(8, 14).Value = 6.933894
$ws.Cells.Item(8, 15).Value = 0.8122763614007964
$ws.Cells.Item(8, 16).Value = 0.863868705896424
$ws.Cells.Item(8, 17).Value = 138.90852442742
$ws.Cells.Item(8, 18).Value = 1250.17671984678
$ws.Cells.Item(8, 19).Value = 0.02647745650484381
$ws.Cells.Item(8, 20).Value = 0.02826380130693834

# Row 9
$ws.Cells.Item(9, 1).Value = "M1"
$ws.Cells.Item(9, 2).Value = "Fn1"
$ws.Cells.Item(9, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(9, 4).Value = "M1"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 60.09979000000001
$ws.Cells.Item(9, 8).Value = 180.29937
$ws.Cells.Item(9, 9).Value = 0.03259661091107292
$ws.Cells.Item(9, 10).Value = 0.03271770480169137
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.02434666666666667
$ws.Cells.Item(9, 14).Value = 0.07304
$ws.Cells.Item(9, 15).Value = 0.008556327142686946
$ws.Cells.Item(9, 16).Value = 0.009099788701510982
$ws.Cells.Item(9, 17).Value = 1.463229553866667
$ws.Cells.Item(9, 18).Value = 13.1690659848
$ws.Cells.Item(9, 19).Value = 0.0002789072666980187
$ws.Cells.Item(9, 20).Value = 0.0002977242004938027

# Row 10
$ws.Cells.Item(10, 1).Value = "M1"
$ws.Cells.Item(10, 2).Value = "Fn1"
$ws.Cells.Item(10, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 60.09979000000001
$ws.Cells.Item(10, 8).Value = 180.29937
$ws.Cells.Item(10, 9).Value = 0.03259661091107292
$ws.Cells.Item(10, 10).Value = 0.03271770480169137
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.509813
$ws.Cells.Item(10, 14).Value = 1.019626
$ws.Cells.Item(10, 15).Value = 0.1791673114565167
$ws.Cells.Item(10, 16).Value = 0.1270315054020651
$ws.Cells.Item(10, 17).Value = 30.63965423927
$ws.Cells.Item(10, 18).Value = 183.83792543562
$ws.Cells.Item(10, 19).Value = 0.005840247139531091
$ws.Cells.Item(10, 20).Value = 0.004156179294259229

# Row 11
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Fn1"
$ws.Cells.Item(11, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 280.168911
$ws.Cells.Item(11, 8).Value = 840.506733
$ws.Cells.Item(11, 9).Value = 0.1519565539454633
$ws.Cells.Item(11, 10).Value = 0.152521060800867
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 2.311298
$ws.Cells.Item(11, 14).Value = 6.933894
$ws.Cells.Item(11, 15).Value = 0.8122763614007964
$ws.Cells.Item(11, 16).Value = 0.863868705896424
$ws.Cells.Item(11, 17).Value = 647.5538436564782
$ws.Cells.Item(11, 18).Value = 5827.984592908303
$ws.Cells.Item(11, 19).Value = 0.1234307167298248
$ws.Cells.Item(11, 20).Value = 0.1317581714159948

# Row 12
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Fn1"
$ws.Cells.Item(12, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(12, 4).Value = "M1"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 280.168911
$ws.Cells.Item(12, 8).Value = 840.506733
$ws.Cells.Item(12, 9).Value = 0.1519565539454633
$ws.Cells.Item(12, 10).Value = 0.152521060800867
$ws.Cells.Item(12, 11).Value = 1
$ws.Cells.Item(12, 12).Value = 0.3333333333333333
$ws.Cells.Item(12, 13).Value = 0.02434666666666667
$ws.Cells.Item(12, 14).Value = 0.07304
$ws.Cells.Item(12, 15).Value = 0.008556327142686946
$ws.Cells.Item(12, 16).Value = 0.009099788701510982
$ws.Cells.Item(12, 17).Value = 6.821179086480001
$ws.Cells.Item(12, 18).Value = 61.39061177832
$ws.Cells.Item(12, 19).Value = 0.001300189987032741
$ws.Cells.Item(12, 20).Value = 0.001387909425818199

# Row 13
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Fn1"
$ws.Cells.Item(13, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 280.168911
$ws.Cells.Item(13, 8).Value = 840.506733
$ws.Cells.Item(13, 9).Value = 0.1519565539454633
$ws.Cells.Item(13, 10).Value = 0.152521060800867
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.509813
$ws.Cells.Item(13, 14).Value = 1.019626
$ws.Cells.Item(13, 15).Value = 0.1791673114565167
$ws.Cells.Item(13, 16).Value = 0.1270315054020651
$ws.Cells.Item(13, 17).Value = 142.833753023643
$ws.Cells.Item(13, 18).Value = 857.002518141858
$ws.Cells.Item(13, 19).Value = 0.02722564722860581
$ws.Cells.Item(13, 20).Value = 0.01937497995905405

# Row 14
$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Fn1"
$ws.Cells.Item(14, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(14, 4).Value = "FAPs"
$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 20.472042
$ws.Cells.Item(14, 8).Value = 40.944084
$ws.Cells.Item(14, 9).Value = 0.01110351945704208
$ws.Cells.Item(14, 10).Value = 0.00742984544919738
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 2.311298
$ws.Cells.Item(14, 14).Value = 6.933894
$ws.Cells.Item(14, 15).Value = 0.8122763614007964
$ws.Cells.Item(14, 16).Value = 0.863868705896424
$ws.Cells.Item(14, 17).Value = 47.31698973051601
$ws.Cells.Item(14, 18).Value = 283.901938383096
$ws.Cells.Item(14, 19).Value = 0.009019126383309088
$ws.Cells.Item(14, 20).Value = 0.006418410973208576

# Row 15
$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Fn1"
$ws.Cells.Item(15, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(15, 4).Value = "M1"
$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 20.472042
$ws.Cells.Item(15, 8).Value = 40.944084
$ws.Cells.Item(15, 9).Value = 0.01110351945704208
$ws.Cells.Item(15, 10).Value = 0.00742984544919738
$ws.Cells.Item(15, 11).Value = 1
$ws.Cells.Item(15, 12).Value = 0.3333333333333333
$ws.Cells.Item(15, 13).Value = 0.02434666666666667
$ws.Cells.Item(15, 14).Value = 0.07304
$ws.Cells.Item(15, 15).Value = 0.008556327142686946
$ws.Cells.Item(15, 16).Value = 0.009099788701510982
$ws.Cells.Item(15, 17).Value = 0.49842598256
$ws.Cells.Item(15, 18).Value = 2.99055589536
$ws.Cells.Item(15, 19).Value = 0.00009500534490964178
$ws.Cells.Item(15, 20).Value = 0.0000676100236725791

# Row 16
$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Fn1"
$ws.Cells.Item(16, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(16, 4).Value = "sCs"
$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 20.472042
$ws.Cells.Item(16, 8).Value = 40.944084
$ws.Cells.Item(16, 9).Value = 0.01110351945704208
$ws.Cells.Item(16, 10).Value = 0.00742984544919738
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 0.509813
$ws.Cells.Item(16, 14).Value = 1.019626
$ws.Cells.Item(16, 15).Value = 0.1791673114565167
$ws.Cells.Item(16, 16).Value = 0.1270315054020651
$ws.Cells.Item(16, 17).Value = 10.436913148146
$ws.Cells.Item(16, 18).Value = 41.747652592584
$ws.Cells.Item(16, 19).Value = 0.001989387728823352
$ws.Cells.Item(16, 20).Value = 0.0009438244523162261

